$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 81
$ws1.Range("F3").Value = 3969
$ws1.Range("F4").Value = 2331
$ws1.Range("F8").Value = 13
$ws1.Range("F9").Value = 193
$ws1.Range("F11").Value = 44
$ws1.Range("F13").Value = 1474
$ws1.Range("F14").Value = 261
$ws1.Range("F15").Value = 2725
$ws1.Range("F16").Value = 187

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 38

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 81
$ws4.Range("F3").Value = 3969
$ws4.Range("F4").Value = 2331
$ws4.Range("F8").Value = 13
$ws4.Range("F9").Value = 38
$ws4.Range("F10").Value = 193
$ws4.Range("F12").Value = 44
$ws4.Range("F16").Value = 1474
$ws4.Range("F17").Value = 261
$ws4.Range("F18").Value = 2725
$ws4.Range("F19").Value = 187
